$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dialogue cell text (adds [set-scroll-delay=...] tags for variable text-scroll speed) ---
$ws.Range("A1").Value = "START_SCENE"
$ws.Range("B1").Value = "CUTSCENE"
$ws.Range("C1").Value = "_"

$ws.Range("A2").Value = "DIALOGUE"
$ws.Range("B2").Value = "Player"
$ws.Range("C2").Value = "[set-scroll-delay=0.05]. . ."
$ws.Range("D2").Value = "_"
$ws.Range("E2").Value = "END_DIALOGUE"

$ws.Range("A3").Value = "DIALOGUE"
$ws.Range("C3").Value = "[block=t][set-scroll-delay=0]click[set-scroll-delay=0.3] . . .[pause=1][set-scroll-delay=0] click.[pause=0.5] click.[pause=0.5][set-scroll-delay=0.05][block=f]"
$ws.Range("D3").Value = "_"
$ws.Range("E3").Value = "END_DIALOGUE"

$ws.Range("A4").Value = "DIALOGUE"
$ws.Range("B4").Value = "Player"
$ws.Range("C4").Value = "Uh . . .[pause=1] Is this thing on?"
$ws.Range("D4").Value = "_"
$ws.Range("E4").Value = "END_DIALOGUE"

$ws.Range("A5").Value = "DIALOGUE"
$ws.Range("B5").Value = "Player"
$ws.Range("C5").Value = "( . . . Crap. I just got this thing too.)"
$ws.Range("D5").Value = "_"
$ws.Range("E5").Value = "END_DIALOGUE"

$ws.Range("A6").Value = "DIALOGUE"
$ws.Range("B6").Value = "Player"
$ws.Range("C6").Value = "(Maybe if I just . . . )"
$ws.Range("D6").Value = "_"
$ws.Range("E6").Value = "END_DIALOGUE"

$ws.Range("A7").Value = "DIALOGUE"
$ws.Range("C7").Value = "[block=t][set-scroll-delay=0]<size=48>[play-sfx=take_damage][screen-shake=0.2,0.3]THWACK!![pause=1]     [play-sfx=take_damage][screen-shake=0.2,0.3]WHAM!!</size>[set-scroll-delay=0.05][block=f]"
$ws.Range("D7").Value = "_"
$ws.Range("E7").Value = "END_DIALOGUE"

$ws.Range("A8").Value = "DIALOGUE"
$ws.Range("B8").Value = "Player"
$ws.Range("C8").Value = "[set-scroll-delay=0.3]. . . [set-scroll-delay=0.05] Did-[pause=0.5]Did it wo[next]"
$ws.Range("D8").Value = "_"
$ws.Range("E8").Value = "END_DIALOGUE"

$ws.Range("A9").Value = "DIALOGUE"
$ws.Range("B9").Value = "???"
$ws.Range("C9").Value = "<size=48>[screen-shake=0.3,0.4]HELLO!</size>"
$ws.Range("D9").Value = "_"
$ws.Range("E9").Value = "END_DIALOGUE"

$ws.Range("A10").Value = "DIALOGUE"
$ws.Range("B10").Value = "???"
$ws.Range("C10").Value = "Welcome back (Uncle name here)! You last activated this <color=blue>TYPOCRYPHA</color> unit 1 year, 1 month, and 13 days ago. You have <color=yellow>1023</color> new updates and <color=yellow>33,333</color> new emails."
$ws.Range("D10").Value = "Naive"
$ws.Range("E10").Value = "END_DIALOGUE"

$ws.Range("A11").Value = "DIALOGUE"
$ws.Range("B11").Value = "Player"
$ws.Range("C11").Value = "O-Oh. No, Im not [Uncles name]. I-uh . . . Dammit! How do I reset-"
$ws.Range("D11").Value = "_"
$ws.Range("E11").Value = "END_DIALOGUE"

$ws.Range("A12").Value = "DIALOGUE"
$ws.Range("B12").Value = "???"
$ws.Range("C12").Value = "<size=48><color=red>ERROR UNIDENTIFIED USER. FORCE SHUTTING DOWN . . .</color></size>"
$ws.Range("D12").Value = "STOP"
$ws.Range("E12").Value = "END_DIALOGUE"

$ws.Range("A13").Value = "DIALOGUE"
$ws.Range("B13").Value = "Player"
$ws.Range("C13").Value = "[block=t]NonononoNONON[fade=out,0,0,0,0][pause=3][next]"
$ws.Range("E13").Value = "END_DIALOGUE"

$ws.Range("A14").Value = "DIALOGUE"
$ws.Range("B14").Value = "???"
$ws.Range("C14").Value = "[fade=in,1.5,0,0,0]REINITIALIZING . . . [pause=3][next]"
$ws.Range("D14").Value = "_"
$ws.Range("E14").Value = "END_DIALOGUE"

$ws.Range("A15").Value = "DIALOGUE"
$ws.Range("B15").Value = "Player"
$ws.Range("C15").Value = "[block=f]OK . . . [pause=1]I think I fixed it. "
$ws.Range("D15").Value = "_"
$ws.Range("E15").Value = "END_DIALOGUE"

$ws.Range("A16").Value = "DIALOGUE"
$ws.Range("B16").Value = "Player"
$ws.Range("C16").Value = ". . . Maybe?"
$ws.Range("D16").Value = "_"
$ws.Range("E16").Value = "END_DIALOGUE"

$ws.Range("A17").Value = "DIALOGUE"
$ws.Range("B17").Value = "Player"
$ws.Range("C17").Value = "Technology sucks."
$ws.Range("D17").Value = "_"
$ws.Range("E17").Value = "END_DIALOGUE"

$ws.Range("A18").Value = "DIALOGUE"
$ws.Range("B18").Value = "Player"
$ws.Range("C18").Value = "[set-scroll-delay=0.3]. . . [set-scroll-delay=0.05] Uh,[pause=0.5] hel[next]"
$ws.Range("D18").Value = "_"
$ws.Range("E18").Value = "END_DIALOGUE"

$ws.Range("A19").Value = "DIALOGUE"
$ws.Range("B19").Value = "???"
$ws.Range("C19").Value = "<size=48>[screen-shake=0.3,0.4]HELLO!</size>"
$ws.Range("D19").Value = "_"
$ws.Range("E19").Value = "END_DIALOGUE"

$ws.Range("A20").Value = "DIALOGUE"
$ws.Range("B20").Value = "???"
$ws.Range("C20").Value = "Thank you for purchasing this <color=blue>TYPOCRYPHA</color> unit. I will be your AI operating assistant."
$ws.Range("D20").Value = "Naive"
$ws.Range("E20").Value = "END_DIALOGUE"

$ws.Range("A21").Value = "DIALOGUE"
$ws.Range("B21").Value = "Clarke"
$ws.Range("C21").Value = "You can call me Clarke!"
$ws.Range("D21").Value = "_"
$ws.Range("E21").Value = "END_DIALOGUE"

$ws.Range("A22").Value = "DIALOGUE"
$ws.Range("B22").Value = "Clarke"
$ws.Range("C22").Value = "What's your name?"
$ws.Range("D22").Value = "_"
$ws.Range("E22").Value = "END_DIALOGUE"

$ws.Range("A23").Value = "END_DIALOGUE"

# --- Clear cells that no longer hold a value in the new layout ---
$ws.Range("B7").ClearContents()
$ws.Range("D13").ClearContents()

# --- Row-height: row 13 reverts to the default (auto) height, row 15 now gets the custom 13.8pt height ---
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(15).RowHeight = 13.8

# --- Move the active selection to C3, matching the saved view state ---
$ws.Range("C3").Select()

